$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Account")

# New "Address" block appended below the existing Account key/value rows (rows 8-16)
$ws.Range("A8").Value = "Address"
$ws.Range("B8").Value = "1 Market St"

$ws.Range("A9").Value = "AddressType"
$ws.Range("A10").Value = "City"
$ws.Range("A11").Value = "Country"
$ws.Range("A12").Value = "Establishment"
$ws.Range("A13").Value = "Unit"
$ws.Range("A14").Value = "Instruction"

$ws.Range("B9").Value = "HOME"
$ws.Range("B10").Value = "San Francisco"
$ws.Range("B11").Value = "United States"
$ws.Range("B12").Value = 123
$ws.Range("B13").Value = 2
$ws.Range("B14").Value = "gate code 4235"

$ws.Range("A15").Value = "Latitude"
$ws.Range("A16").Value = "Longitude"

# Selection on the Account sheet moves to A17 (just below the new data)
$ws.Range("A17").Select() | Out-Null

# The Patient tab becomes the active/selected tab, with selection at N4
$ws2 = $wb.Worksheets.Item("Patient")
$ws2.Range("N4").Select() | Out-Null
$ws2.Activate()
